# Auto-generated edit script: applies numeric corrections to the
# currentAveragePrice / LevePrice / LeveProfit columns (H:N) across
# several rows on multiple craft sheets, per the scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 47
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 15000
$ws.Range("L47").Value = 15000
$ws.Range("N47").Value = -16944

# Row 64
$ws.Range("H64").Value = 2975

# Row 67
$ws.Range("H67").Value = 2975

# Row 132
$ws.Range("H132").Value = 3974.76
$ws.Range("I132").Value = 3835
$ws.Range("K132").Value = 11505
$ws.Range("M132").Value = -8975

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1971671.2
$ws.Range("I32").Value = 2423644.5
$ws.Range("K32").Value = 2423644.5
$ws.Range("M32").Value = -2423357.5

# Row 45
$ws.Range("H45").Value = 1808.5
$ws.Range("I45").Value = 1342.2142
$ws.Range("J45").Value = 2896.5
$ws.Range("K45").Value = 1342.2142
$ws.Range("L45").Value = 2896.5
$ws.Range("M45").Value = -965.2141999999999
$ws.Range("N45").Value = -3650.5

# Row 61
$ws.Range("H61").Value = 6194.737
$ws.Range("I61").Value = 4366.6665
$ws.Range("J61").Value = 6537.5
$ws.Range("K61").Value = 4366.6665
$ws.Range("L61").Value = 6537.5
$ws.Range("M61").Value = -4154.6665
$ws.Range("N61").Value = -6961.5

# Row 74
$ws.Range("H74").Value = 2991.3076
$ws.Range("I74").Value = 2101.8823
$ws.Range("J74").Value = 4671.3335
$ws.Range("K74").Value = 2101.8823
$ws.Range("L74").Value = 4671.3335
$ws.Range("M74").Value = -1227.8823
$ws.Range("N74").Value = -6419.3335

# Row 77
$ws.Range("H77").Value = 2991.3076
$ws.Range("I77").Value = 2101.8823
$ws.Range("J77").Value = 4671.3335
$ws.Range("K77").Value = 10509.4115
$ws.Range("L77").Value = 23356.6675
$ws.Range("M77").Value = -6141.411500000002
$ws.Range("N77").Value = -32092.6675

# Row 122
$ws.Range("H122").Value = 49110.19
$ws.Range("I122").Value = 56973
$ws.Range("J122").Value = 1933.3334
$ws.Range("K122").Value = 170919
$ws.Range("L122").Value = 5800.0002
$ws.Range("M122").Value = -168469
$ws.Range("N122").Value = -10700.0002

# Row 132
$ws.Range("H132").Value = 4029.0278
$ws.Range("I132").Value = 3301.32
$ws.Range("K132").Value = 9903.960000000001
$ws.Range("M132").Value = -7373.960000000001

# Row 136
$ws.Range("H136").Value = 6194.737
$ws.Range("I136").Value = 4366.6665
$ws.Range("J136").Value = 6537.5
$ws.Range("K136").Value = 13099.9995
$ws.Range("L136").Value = 19612.5
$ws.Range("M136").Value = -10549.9995
$ws.Range("N136").Value = -24712.5

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3112.8438
$ws.Range("I134").Value = 3070.7778
$ws.Range("K134").Value = 9212.3334
$ws.Range("M134").Value = -6677.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 1751.6666
$ws.Range("J122").Value = 1914.3125
$ws.Range("L122").Value = 5742.9375
$ws.Range("N122").Value = -10642.9375

# Row 132
$ws.Range("H132").Value = 5557974
$ws.Range("I132").Value = 2283.3125
$ws.Range("J132").Value = 11907334
$ws.Range("K132").Value = 6849.9375
$ws.Range("L132").Value = 35722002
$ws.Range("M132").Value = -4319.9375
$ws.Range("N132").Value = -35727062

# Row 134
$ws.Range("H134").Value = 6861.0415
$ws.Range("I134").Value = 9009.714
$ws.Range("J134").Value = 3852.9
$ws.Range("K134").Value = 27029.142
$ws.Range("L134").Value = 11558.7
$ws.Range("M134").Value = -24494.142
$ws.Range("N134").Value = -16628.7

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 8756.904
$ws.Range("J3").Value = 9791
$ws.Range("L3").Value = 29373
$ws.Range("N3").Value = -29597

# Row 109
$ws.Range("H109").Value = 2125.0527
$ws.Range("I109").Value = 904.3333
$ws.Range("J109").Value = 2688.4614
$ws.Range("K109").Value = 2712.9999
$ws.Range("L109").Value = 8065.3842
$ws.Range("M109").Value = -1672.9999
$ws.Range("N109").Value = -10145.3842

# Row 137
$ws.Range("H137").Value = 31686.975
$ws.Range("I137").Value = 6671.8096
$ws.Range("J137").Value = 62588.06
$ws.Range("K137").Value = 20015.4288
$ws.Range("L137").Value = 187764.18
$ws.Range("M137").Value = -14915.4288
$ws.Range("N137").Value = -197964.18

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 80402.36
$ws.Range("I113").Value = 101857.55
$ws.Range("J113").Value = 1733.3334
$ws.Range("K113").Value = 101857.55
$ws.Range("L113").Value = 1733.3334
$ws.Range("M113").Value = -99687.55
$ws.Range("N113").Value = -6073.3334

# Row 122
$ws.Range("H122").Value = 9615.923000000001
$ws.Range("I122").Value = 14300
$ws.Range("J122").Value = 2121.4
$ws.Range("K122").Value = 42900
$ws.Range("L122").Value = 6364.200000000001
$ws.Range("M122").Value = -40450
$ws.Range("N122").Value = -11264.2

# Row 132
$ws.Range("H132").Value = 2759.8147
$ws.Range("I132").Value = 2361.1667
$ws.Range("K132").Value = 7083.500100000001
$ws.Range("M132").Value = -4553.500100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8249
$ws.Range("I7").Value = 8498.666999999999
$ws.Range("J7").Value = 7500
$ws.Range("K7").Value = 8498.666999999999
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = -8386.666999999999
$ws.Range("N7").Value = -7724

# Row 126
$ws.Range("H126").Value = 8249
$ws.Range("I126").Value = 8498.666999999999
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 25496.001
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -23026.001
$ws.Range("N126").Value = -27440

# Row 132
$ws.Range("H132").Value = 3261.9395
$ws.Range("I132").Value = 2799.4285
$ws.Range("J132").Value = 4071.3333
$ws.Range("K132").Value = 8398.2855
$ws.Range("L132").Value = 12213.9999
$ws.Range("M132").Value = -5868.2855
$ws.Range("N132").Value = -17273.9999

# Row 136
$ws.Range("H136").Value = 4763517
$ws.Range("I136").Value = 1208.6154
$ws.Range("K136").Value = 3625.8462
$ws.Range("M136").Value = -1075.8462

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3898
$ws.Range("I81").Value = 4924.2856
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 9848.5712
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -8787.5712
$ws.Range("N81").Value = -8122

# Row 84
$ws.Range("H84").Value = 3898
$ws.Range("I84").Value = 4924.2856
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 49242.856
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -43938.856
$ws.Range("N84").Value = -40608

# Row 123
$ws.Range("H123").Value = 41204.145
$ws.Range("J123").Value = 41204.145
$ws.Range("L123").Value = 41204.145
$ws.Range("N123").Value = -51004.145

# Row 126
$ws.Range("H126").Value = 1037.4445
$ws.Range("I126").Value = 1037.4445
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3112.3335
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -642.3335000000002
$ws.Range("N126").ClearContents()

# Row 132
$ws.Range("H132").Value = 3971241.5
$ws.Range("I132").Value = 3157.2632
$ws.Range("J132").Value = 7249224
$ws.Range("K132").Value = 9471.7896
$ws.Range("L132").Value = 21747672
$ws.Range("M132").Value = -6941.7896
$ws.Range("N132").Value = -21752732

# Row 136
$ws.Range("H136").Value = 2848.4905
$ws.Range("I136").Value = 2688.2812
$ws.Range("J136").Value = 3092.6191
$ws.Range("K136").Value = 8064.8436
$ws.Range("L136").Value = 9277.8573
$ws.Range("M136").Value = -5514.8436
$ws.Range("N136").Value = -14377.8573
